$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.323.40"
$ws.Range("E2").Value = "'  -7.23%  "
$ws.Range("D3").Value = "'2.884.87"
$ws.Range("E3").Value = "'  -5.51%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'551.99"
$ws.Range("E5").Value = "'  -5.80%  "
$ws.Range("D6").Value = "'121.18"
$ws.Range("E6").Value = "'  -6.91%  "
$ws.Range("E7").Value = "'  +0.19%  "
$ws.Range("D8").Value = "'2.874.23"
$ws.Range("E8").Value = "'  -5.82%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "'  -2.73%  "
$ws.Range("E10").Value = "'  -10.79%  "
$ws.Range("E11").Value = "'  -9.98%  "
$ws.Range("D12").Value = "'0.432"
$ws.Range("E12").Value = "'  -2.06%  "
$ws.Range("E13").Value = "'  -10.78%  "
$ws.Range("D14").Value = "'31.31"
$ws.Range("E14").Value = "'  -7.22%  "
$ws.Range("E15").Value = "'  -0.93%  "
$ws.Range("D16").Value = "'3.352.02"
$ws.Range("E16").Value = "'  -5.78%  "
$ws.Range("D17").Value = "'2.876.70"
$ws.Range("E17").Value = "'  -5.88%  "
$ws.Range("D18").Value = "'57.259.53"
$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "'  +0.60%  "
$ws.Range("D20").Value = "'408.03"
$ws.Range("E20").Value = "'  -9.18%  "
$ws.Range("D21").Value = "'12.76"
$ws.Range("E21").Value = "'  -5.79%  "
$ws.Range("D22").Value = "'0.650"
$ws.Range("E22").Value = "'  -3.57%  "
$ws.Range("D23").Value = "'6.71"
$ws.Range("E23").Value = "'  -8.80%  "
$ws.Range("D24").Value = "'12.54"
$ws.Range("E24").Value = "'  -2.72%  "
$ws.Range("D25").Value = "'76.59"
$ws.Range("E25").Value = "'  -5.60%  "
$ws.Range("E26").Value = "'  -0.32%  "
$ws.Range("E27").Value = "'  -0.04%  "
$ws.Range("E28").Value = "'  -4.44%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'1.90"
$ws.Range("E29").Value = "'  -5.86%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.11"
$ws.Range("E30").Value = "'  -4.23%  "
$ws.Range("D31").Value = "'6.04"
$ws.Range("E31").Value = "'  -6.62%  "
$ws.Range("D32").Value = "'24.56"
$ws.Range("E32").Value = "'  -5.34%  "
$ws.Range("D33").Value = "'0.0947"
$ws.Range("E33").Value = "'  -2.98%  "
$ws.Range("D34").Value = "'2.01"
$ws.Range("E34").Value = "'  -14.04%  "
$ws.Range("E35").Value = "'  -6.72%  "
$ws.Range("D36").Value = "'0.892"
$ws.Range("E36").Value = "'  -8.60%  "
$ws.Range("D37").Value = "'48.32"
$ws.Range("E37").Value = "'  -4.09%  "
$ws.Range("E38").Value = "'  +5.04%  "
$ws.Range("D39").Value = "'0.0₃0610"
$ws.Range("E39").Value = "'  -12.28%  "
$ws.Range("D40").Value = "'0.0343"
$ws.Range("E40").Value = "'  -8.93%  "
$ws.Range("E41").Value = "'  -3.31%  "
$ws.Range("D42").Value = "'2.588.46"
$ws.Range("E42").Value = "'  -4.14%  "
$ws.Range("D43").Value = "'357.82"
$ws.Range("E43").Value = "'  -6.47%  "
$ws.Range("E45").Value = "'  -7.68%  "
$ws.Range("D46").Value = "'117.18"
$ws.Range("E46").Value = "'  -5.63%  "
$ws.Range("E47").Value = "'  -5.51%  "
$ws.Range("E48").Value = "'  -2.85%  "
$ws.Range("D49").Value = "'1.91"
$ws.Range("E49").Value = "'  -5.20%  "
$ws.Range("D50").Value = "'22.33"
$ws.Range("E50").Value = "'  -7.17%  "
$ws.Range("D51").Value = "'1.93"
$ws.Range("E51").Value = "'  -7.43%  "
